# Refresh the crypto market snapshot on Sheet1 with the latest scrape.
# Columns B (coin) and C (link) occasionally shift rows when the ranking
# reorders; columns D (price) and E (1h change) are updated for every row.
#
# Some new Price strings are plain numeric text (e.g. "1.00", "0.0779")
# that Excel would otherwise auto-coerce into a real number on assignment
# (dropping the trailing zero / flipping to scientific notation). Those are
# written with a leading apostrophe to force text, then the cell style is
# reset to "Normal" so the quote-prefix marker doesn't leave a stray style
# index behind (matching the source cells, which carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Val='68.403.10'; Text=$false},
    @{Cell='E2'; Val='  +0.04%  '; Text=$false},
    @{Cell='D3'; Val='2.654.15'; Text=$false},
    @{Cell='E3'; Val='  +0.48%  '; Text=$false},
    @{Cell='D4'; Val='''0.999'; Text=$true},
    @{Cell='E4'; Val='  -0.06%  '; Text=$false},
    @{Cell='D5'; Val='''598.14'; Text=$true},
    @{Cell='E5'; Val='  -0.37%  '; Text=$false},
    @{Cell='D6'; Val='''159.38'; Text=$true},
    @{Cell='E6'; Val='  +3.01%  '; Text=$false},
    @{Cell='E8'; Val='  -0.35%  '; Text=$false},
    @{Cell='D9'; Val='''0.146'; Text=$true},
    @{Cell='E9'; Val='  +5.78%  '; Text=$false},
    @{Cell='E10'; Val='  -1.10%  '; Text=$false},
    @{Cell='D11'; Val='''5.27'; Text=$true},
    @{Cell='E11'; Val='  +0.52%  '; Text=$false},
    @{Cell='E12'; Val='  +0.60%  '; Text=$false},
    @{Cell='D13'; Val='''28.15'; Text=$true},
    @{Cell='E13'; Val='  +0.28%  '; Text=$false},
    @{Cell='D14'; Val='''0.0000191'; Text=$true},
    @{Cell='E14'; Val='  +1.42%  '; Text=$false},
    @{Cell='D15'; Val='3.135.92'; Text=$false},
    @{Cell='E15'; Val='  +0.39%  '; Text=$false},
    @{Cell='D16'; Val='68.331.67'; Text=$false},
    @{Cell='E16'; Val='  +0.22%  '; Text=$false},
    @{Cell='D17'; Val='2.649.03'; Text=$false},
    @{Cell='E17'; Val='  +0.44%  '; Text=$false},
    @{Cell='D18'; Val='''11.44'; Text=$true},
    @{Cell='E18'; Val='  -0.17%  '; Text=$false},
    @{Cell='D19'; Val='''365.16'; Text=$true},
    @{Cell='E19'; Val='  -0.56%  '; Text=$false},
    @{Cell='D20'; Val='''7.36'; Text=$true},
    @{Cell='E20'; Val='  -1.02%  '; Text=$false},
    @{Cell='E21'; Val='  +3.40%  '; Text=$false},
    @{Cell='D22'; Val='''4.84'; Text=$true},
    @{Cell='E22'; Val='  -0.47%  '; Text=$false},
    @{Cell='E23'; Val='  -2.31%  '; Text=$false},
    @{Cell='D24'; Val='''75.18'; Text=$true},
    @{Cell='E24'; Val='  +2.17%  '; Text=$false},
    @{Cell='E25'; Val='  +0.03%  '; Text=$false},
    @{Cell='D26'; Val='''9.77'; Text=$true},
    @{Cell='E26'; Val='  -2.60%  '; Text=$false},
    @{Cell='B27'; Val='PEPE'; Text=$false},
    @{Cell='C27'; Val='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; Text=$false},
    @{Cell='D27'; Val='''0.0000105'; Text=$true},
    @{Cell='E27'; Val='  +0.81%  '; Text=$false},
    @{Cell='B28'; Val='WrappedeETH'; Text=$false},
    @{Cell='C28'; Val='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; Text=$false},
    @{Cell='D28'; Val='2.784.32'; Text=$false},
    @{Cell='E28'; Val='  +0.51%  '; Text=$false},
    @{Cell='D29'; Val='''1.02'; Text=$true},
    @{Cell='E29'; Val='  +1.88%  '; Text=$false},
    @{Cell='D30'; Val='''560.46'; Text=$true},
    @{Cell='E30'; Val='  -2.30%  '; Text=$false},
    @{Cell='D31'; Val='''8.05'; Text=$true},
    @{Cell='E31'; Val='  +0.43%  '; Text=$false},
    @{Cell='E32'; Val='  -0.84%  '; Text=$false},
    @{Cell='D33'; Val='''1.87'; Text=$true},
    @{Cell='E33'; Val='  +0.34%  '; Text=$false},
    @{Cell='E34'; Val='  -0.94%  '; Text=$false},
    @{Cell='E35'; Val='  -0.05%  '; Text=$false},
    @{Cell='D36'; Val='''1.58'; Text=$true},
    @{Cell='E36'; Val='  +1.87%  '; Text=$false},
    @{Cell='D37'; Val='''19.89'; Text=$true},
    @{Cell='E37'; Val='  +2.92%  '; Text=$false},
    @{Cell='D38'; Val='''159.64'; Text=$true},
    @{Cell='E38'; Val='  -0.52%  '; Text=$false},
    @{Cell='D39'; Val='''0.372'; Text=$true},
    @{Cell='E39'; Val='  +0.79%  '; Text=$false},
    @{Cell='E40'; Val='  -2.22%  '; Text=$false},
    @{Cell='E41'; Val='  -0.66%  '; Text=$false},
    @{Cell='D42'; Val='0.0₆0337'; Text=$false},
    @{Cell='E42'; Val='  +5.08%  '; Text=$false},
    @{Cell='B43'; Val='WhiteBITCoin'; Text=$false},
    @{Cell='C43'; Val='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; Text=$false},
    @{Cell='D43'; Val='''17.81'; Text=$true},
    @{Cell='E43'; Val='  +0.35%  '; Text=$false},
    @{Cell='B44'; Val='dogwifhat'; Text=$false},
    @{Cell='C44'; Val='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; Text=$false},
    @{Cell='D44'; Val='''2.64'; Text=$true},
    @{Cell='E44'; Val='  -0.20%  '; Text=$false},
    @{Cell='B45'; Val='USDe'; Text=$false},
    @{Cell='C45'; Val='https://coinranking.com/coin/exbfr2U-0+usde-usde'; Text=$false},
    @{Cell='D45'; Val='''1.00'; Text=$true},
    @{Cell='E45'; Val='  +0.03%  '; Text=$false},
    @{Cell='B46'; Val='Aave'; Text=$false},
    @{Cell='C46'; Val='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; Text=$false},
    @{Cell='D46'; Val='''158.36'; Text=$true},
    @{Cell='E46'; Val='  +0.00%  '; Text=$false},
    @{Cell='B47'; Val='Filecoin'; Text=$false},
    @{Cell='C47'; Val='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Text=$false},
    @{Cell='D47'; Val='''3.78'; Text=$true},
    @{Cell='E47'; Val='  +0.16%  '; Text=$false},
    @{Cell='B48'; Val='InjectiveProtocol'; Text=$false},
    @{Cell='C48'; Val='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; Text=$false},
    @{Cell='D48'; Val='''22.28'; Text=$true},
    @{Cell='E48'; Val='  +1.31%  '; Text=$false},
    @{Cell='B49'; Val='Optimism'; Text=$false},
    @{Cell='C49'; Val='https://coinranking.com/coin/n1p-s_gm1+optimism-op'; Text=$false},
    @{Cell='D49'; Val='''1.70'; Text=$true},
    @{Cell='E49'; Val='  -1.14%  '; Text=$false},
    @{Cell='B50'; Val='Cronos'; Text=$false},
    @{Cell='C50'; Val='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Text=$false},
    @{Cell='D50'; Val='''0.0779'; Text=$true},
    @{Cell='E50'; Val='  -0.12%  '; Text=$false},
    @{Cell='B51'; Val='Mantle'; Text=$false},
    @{Cell='C51'; Val='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; Text=$false},
    @{Cell='D51'; Val='''0.617'; Text=$true},
    @{Cell='E51'; Val='  +0.13%  '; Text=$false}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = $u.Val
    if ($u.Text) {
        $cell.Style = "Normal"
    }
}
